$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") updates. A handful of the new prices are plain
# decimal-looking strings ("1.00", "6.21", ...); Excel would silently
# reinterpret a bare Value assignment like that as a number, so those
# are written with a leading apostrophe (forces text) and then the
# cell style is put back to Normal so no stray number-format/quote-
# prefix formatting is left behind - only the text value changes.
$ws.Range("D2").Value = "65.054.45"
$ws.Range("D3").Value = "3.179.71"
$ws.Range("D5").Value = "'578.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'151.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "3.178.91"
$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'6.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'38.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.699.49"
$ws.Range("D16").Value = "65.167.07"
$ws.Range("D17").Value = "3.183.41"
$ws.Range("D18").Value = "'7.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = "'514.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'14.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'15.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'85.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Value = "'28.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Value = "'6.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Value = "'55.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Value = "'478.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Value = "'8.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "3.063.05"
$ws.Range("D43").Value = "'0.119"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.289"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'2.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'29.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "0.0₃0612"
$ws.Range("D51").Value = "'120.79"
$ws.Range("D51").Style = "Normal"

# Column E ("Volume(1h)") updates - these are always padded with spaces
# and a percent sign, so Excel keeps them as plain text automatically.
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +3.83%  "
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("E10").Value = "  +5.48%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("E14").Value = "  +6.27%  "
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("E17").Value = "  +3.76%  "
$ws.Range("E18").Value = "  +5.42%  "
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E20").Value = "  +7.44%  "
$ws.Range("E21").Value = "  +6.63%  "
$ws.Range("E22").Value = "  +7.31%  "
$ws.Range("E23").Value = "  +6.66%  "
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("E27").Value = "  +11.64%  "
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("E29").Value = "  +7.66%  "
$ws.Range("E30").Value = "  +6.62%  "
$ws.Range("E31").Value = "  +13.64%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +6.06%  "
$ws.Range("E34").Value = "  +8.30%  "
$ws.Range("E35").Value = "  +6.16%  "
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  +10.50%  "
$ws.Range("E38").Value = "  +5.64%  "
$ws.Range("E39").Value = "  +11.70%  "
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("E41").Value = "  +4.49%  "
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("E44").Value = "  +8.22%  "
$ws.Range("E45").Value = "  +7.52%  "
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("E47").Value = "  +18.12%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  +8.78%  "
$ws.Range("E51").Value = "  +1.43%  "
